# Crit bonus changed slightly
# +1 damage when phys, special effect otherwise
#
# The "Class Tree" labels under each class snippet are widened/centered so
# the (now longer) crit-bonus description fits, and the second "Notes"
# caption (the smaller Corbel-font one) is shrunk from 10.5pt to 8pt.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Class Tree" textbox (left column) ---------------------------------
$classTree1 = $s.Shapes.Item("TextBox 22")
$classTree1.Left = 26.57503987007874   # -> a:off x = 337503 EMU
$classTree1.Width = 173.92496492992126 # -> a:ext cx = 2208847 EMU
$classTree1.TextFrame.WordWrap = -1    # wrap="none" -> wrap="square"
$classTree1.TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter -> algn="ctr"

# --- "Class Tree" textbox (right column) ---------------------------------
$classTree2 = $s.Shapes.Item("TextBox 25")
$classTree2.Left = 226.3931504062992   # -> a:off x = 2875193 EMU
$classTree2.Width = 173.92496492992126 # -> a:ext cx = 2208847 EMU
$classTree2.TextFrame.WordWrap = -1    # wrap="none" -> wrap="square"
$classTree2.TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter -> algn="ctr"

# --- "Notes" caption (second, smaller "Notes" run in Rectangle 33) ------
$notesBox = $s.Shapes.Item("Rectangle 33")
$notesRange = $notesBox.TextFrame.TextRange
# Full text is "NotesNotes" (para1 "Notes", para2 empty, para3 "Notes");
# the second "Notes" run sits at characters 8-12.
$secondNotes = $notesRange.Characters(8, 5)
$secondNotes.Font.Size = 8
